# Bump the "Version" and "Date" metadata values on the Metadata sheet,
# matching the new 1.1.0 release of the KLEvaluationTypeCodes term.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 3: Property "Version" -> Value "1.1.0" (was "1.0.0")
$ws.Range("B3").Value = "1.1.0"

# Row 8: Property "Date" -> Value "2023-07-10T23:08:03+02:00" (was "2023-06-07T11:52:14+02:00")
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
